# Applies the crypto price/volume refresh from the Aug 31 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (so numeric-looking strings like "65.00"
# keep their literal formatting instead of being coerced to a number),
# then restore the cell to the default "Normal" style so no visible
# formatting change is introduced.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '58.896.89'
$ws.Range('E2').Value = '  +0.18%  '

$ws.Range('D3').Value = '2.497.93'
$ws.Range('E3').Value = '  +0.22%  '

$ws.Range('E4').Value = '  -0.04%  '

Set-TextValue $ws.Range('D5') '532.39'
$ws.Range('E5').Value = '  -0.31%  '

Set-TextValue $ws.Range('D6') '134.91'
$ws.Range('E6').Value = '  -1.23%  '

$ws.Range('E7').Value = '  +0.26%  '

$ws.Range('E8').Value = '  +1.07%  '

Set-TextValue $ws.Range('D9') '0.101'
$ws.Range('E9').Value = '  +0.83%  '

$ws.Range('E10').Value = '  -0.97%  '

$ws.Range('E11').Value = '  +1.55%  '

Set-TextValue $ws.Range('D12') '0.347'
$ws.Range('E12').Value = '  +0.30%  '

$ws.Range('D13').Value = '2.939.81'
$ws.Range('E13').Value = '  +0.05%  '

$ws.Range('D14').Value = '58.831.13'
$ws.Range('E14').Value = '  +0.15%  '

$ws.Range('E15').Value = '  -1.61%  '

$ws.Range('E16').Value = '  -0.91%  '

$ws.Range('D17').Value = '2.489.75'
$ws.Range('E17').Value = '  -0.68%  '

$ws.Range('E18').Value = '  +0.24%  '

$ws.Range('E19').Value = '  -0.10%  '

Set-TextValue $ws.Range('D20') '322.51'
$ws.Range('E20').Value = '  -0.61%  '

$ws.Range('E21').Value = '  -0.01%  '

$ws.Range('E22').Value = '  +1.34%  '

Set-TextValue $ws.Range('D23') '65.00'
$ws.Range('E23').Value = '  +2.35%  '

$ws.Range('E24').Value = '  +0.79%  '

Set-TextValue $ws.Range('D25') '0.163'
$ws.Range('E25').Value = '  -0.52%  '

$ws.Range('E26').Value = '  +0.95%  '

$ws.Range('E27').Value = '  -0.76%  '

$ws.Range('E28').Value = '  -1.52%  '

Set-TextValue $ws.Range('D29') '170.16'
$ws.Range('E29').Value = '  +2.36%  '

$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D30') '6.44'
$ws.Range('E30').Value = '  -4.11%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D31') '1.74'
$ws.Range('E31').Value = '  -1.26%  '

$ws.Range('E32').Value = '  +0.94%  '

$ws.Range('E33').Value = '  +0.10%  '

Set-TextValue $ws.Range('D34') '18.32'
$ws.Range('E34').Value = '  -0.66%  '

$ws.Range('E35').Value = '  -2.13%  '

$ws.Range('E36').Value = '  -1.28%  '

Set-TextValue $ws.Range('D37') '1.51'
$ws.Range('E37').Value = '  -2.85%  '

$ws.Range('E38').Value = '  -2.08%  '

$ws.Range('E39').Value = '  -1.36%  '

Set-TextValue $ws.Range('D40') '281.03'
$ws.Range('E40').Value = '  +1.18%  '

$ws.Range('E41').Value = '  +0.42%  '

Set-TextValue $ws.Range('D42') '4.99'
$ws.Range('E42').Value = '  -4.52%  '

Set-TextValue $ws.Range('D43') '129.53'
$ws.Range('E43').Value = '  +2.62%  '

Set-TextValue $ws.Range('D44') '10.89'
$ws.Range('E44').Value = '  +0.47%  '

$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('E46').Value = '  +0.09%  '

$ws.Range('E47').Value = '  -2.46%  '

$ws.Range('E48').Value = '  -1.52%  '

Set-TextValue $ws.Range('D49') '17.20'
$ws.Range('E49').Value = '  -1.03%  '

$ws.Range('D50').Value = '1.751.79'
$ws.Range('E50').Value = '  -0.59%  '

$ws.Range('E51').Value = '  -0.41%  '
